$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "37.114.94"
Set-TextValue $ws.Range("E2") "  -0.63%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.025.40"
Set-TextValue $ws.Range("E3") "  -1.01%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.13%  "

# Row 5
Set-TextValue $ws.Range("D5") "226.55"
Set-TextValue $ws.Range("E5") "  -1.00%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.610"
Set-TextValue $ws.Range("E6") "  -0.87%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.05%  "

# Row 8
Set-TextValue $ws.Range("D8") "55.16"
Set-TextValue $ws.Range("E8") "  -3.12%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.378"
Set-TextValue $ws.Range("E9") "  -2.09%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0785"
Set-TextValue $ws.Range("E10") "  -0.39%  "

# Row 11
Set-TextValue $ws.Range("E11") "  -4.37%  "

# Row 12
Set-TextValue $ws.Range("D12") "2.329.65"
Set-TextValue $ws.Range("E12") "  -0.70%  "

# Row 13
Set-TextValue $ws.Range("D13") "14.15"
Set-TextValue $ws.Range("E13") "  -3.92%  "

# Row 14
Set-TextValue $ws.Range("D14") "20.20"
Set-TextValue $ws.Range("E14") "  -2.10%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.745"
Set-TextValue $ws.Range("E15") "  -1.45%  "

# Row 16
Set-TextValue $ws.Range("D16") "5.18"
Set-TextValue $ws.Range("E16") "  -2.30%  "

# Row 17
Set-TextValue $ws.Range("D17") "1.997.53"
Set-TextValue $ws.Range("E17") "  -2.59%  "

# Row 18
Set-TextValue $ws.Range("D18") "37.076.81"
Set-TextValue $ws.Range("E18") "  -0.33%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.46"
Set-TextValue $ws.Range("E19") "  +6.11%  "

# Row 20
Set-TextValue $ws.Range("D20") "68.81"
Set-TextValue $ws.Range("E20") "  -1.00%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.0₃0816"
Set-TextValue $ws.Range("E21") "  -1.34%  "

# Row 22
Set-TextValue $ws.Range("D22") "223.11"
Set-TextValue $ws.Range("E22") "  -1.21%  "

# Row 23
Set-TextValue $ws.Range("E23") "  -0.03%  "

# Row 24
Set-TextValue $ws.Range("E24") "  +1.91%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.19"
Set-TextValue $ws.Range("E25") "  -4.62%  "

# Row 26
Set-TextValue $ws.Range("D26") "165.27"
Set-TextValue $ws.Range("E26") "  -1.62%  "

# Row 27
Set-TextValue $ws.Range("D27") "9.18"
Set-TextValue $ws.Range("E27") "  -5.00%  "

# Row 28
Set-TextValue $ws.Range("E28") "  -0.84%  "

# Row 29
Set-TextValue $ws.Range("D29") "18.68"
Set-TextValue $ws.Range("E29") "  -1.47%  "

# Row 30
Set-TextValue $ws.Range("E30") "  -3.09%  "

# Row 31
Set-TextValue $ws.Range("E31") "  -1.17%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.51"
Set-TextValue $ws.Range("E32") "  -0.58%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.0606"
Set-TextValue $ws.Range("E33") "  -1.25%  "

# Row 34
Set-TextValue $ws.Range("E34") "  -1.70%  "

# Row 35
Set-TextValue $ws.Range("D35") "2.34"
Set-TextValue $ws.Range("E35") "  -3.62%  "

# Row 36
Set-TextValue $ws.Range("E36") "  +0.79%  "

# Row 37
Set-TextValue $ws.Range("E37") "  +0.32%  "

# Row 38
Set-TextValue $ws.Range("D38") "5.55"
Set-TextValue $ws.Range("E38") "  +5.32%  "

# Row 39
Set-TextValue $ws.Range("D39") "3.10"
Set-TextValue $ws.Range("E39") "  -4.37%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.464.39"
Set-TextValue $ws.Range("E40") "  -0.60%  "

# Row 41
Set-TextValue $ws.Range("E41") "  -3.27%  "

# Row 42
Set-TextValue $ws.Range("D42") "95.41"
Set-TextValue $ws.Range("E42") "  -0.83%  "

# Row 43
Set-TextValue $ws.Range("B43") "HuobiToken"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D43") "2.81"
Set-TextValue $ws.Range("E43") "  -2.70%  "

# Row 44
Set-TextValue $ws.Range("B44") "InjectiveProtocol"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D44") "16.31"
Set-TextValue $ws.Range("E44") "  -4.84%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.0911"
Set-TextValue $ws.Range("E45") "  -3.32%  "

# Row 46
Set-TextValue $ws.Range("E46") "  -1.87%  "

# Row 47
Set-TextValue $ws.Range("D47") "7.25"
Set-TextValue $ws.Range("E47") "  +1.77%  "

# Row 48
Set-TextValue $ws.Range("E48") "  -0.98%  "

# Row 49
Set-TextValue $ws.Range("D49") "2.94"
Set-TextValue $ws.Range("E49") "  +0.67%  "

# Row 50
Set-TextValue $ws.Range("D50") "2.215.34"
Set-TextValue $ws.Range("E50") "  -0.57%  "

# Row 51
Set-TextValue $ws.Range("D51") "3.61"
Set-TextValue $ws.Range("E51") "  -8.18%  "
